$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55. This shifts the existing rows 55-104
# down to 56-105 (and extends the used range / dimension to A1:R105),
# matching every row's data moving down by one position as seen in the
# diff (e.g. old row 55 -> new row 56, ..., old row 104 -> new row 105).
$ws.Rows("55:55").Insert()

# Populate the newly inserted row 55 with the new data point that was
# added to the dataset (same "Puerro" / Vega Central Mapocho record
# shape as every other row, with its own date and price figures).
$ws.Range("A55").Value = 9
$ws.Range("B55").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C55").Value = "Metropolitana"
$ws.Range("D55").Value = 44790
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = 100112005
$ws.Range("G55").Value = "Puerro"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 160
$ws.Range("K55").Value = 7000
$ws.Range("L55").Value = 7000
$ws.Range("M55").Value = 7000
$ws.Range("N55").Value = "`$/paquete 20 unidades"
$ws.Range("O55").Value = "Provincia de Chacabuco"
$ws.Range("P55").Value = 350
$ws.Range("Q55").Value = 20
$ws.Range("R55").Value = "Hortaliza"
